$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: add the comment explaining the health-points feature ---
$ws.Range("E16").Value = "J'ai rajouté le nombre de vies totaux des bateaux après une inspiration de Sebastien Moraz. (Parce que c'est quand même pratique de voir les points de vie)"
$ws.Rows.Item(16).RowHeight = 60

# --- Row 17: fill in the new work-log entry ---
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = "Création des conditions si le chiffre choisi est inférieur à 10 et affichage des points de vie des bateaux"
$ws.Range("C17").WrapText = $true
$ws.Range("D17").Value = "1h30"
$ws.Range("E17").Value = "j'ai essayé de mettre des couleurs aux tâches touchées."
$ws.Range("F17").Value = "Toujours entrain d'essayer de mettre une couleur aux cases touchées"
$ws.Range("F17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 45

# --- Update the active selection like the author left it ---
$ws.Range("F17").Select()
